# Krull_2005.xlsx - add a new "pro_usda_soil_order" controlled-vocabulary
# column to the "profile" sheet, backed by a new column of soil-order
# values on the "controlled vocabulary" sheet, with a data-validation
# dropdown tying the two together.

$wb = $excel.ActiveWorkbook

# --- 1. "profile" sheet: insert a new column N -------------------------
$profile = $wb.Worksheets.Item("profile")
$profile.Columns("N").Insert()
$profile.Range("N1").Value = "pro_usda_soil_order"

# --- 2. "controlled vocabulary" sheet: insert a new column E ----------
$cv = $wb.Worksheets.Item("controlled vocabulary")
$cv.Columns("E").Insert()
$cv.Range("E2").Value = "pro_usda_soil_order"
$cv.Range("E4").Value = "Alfisols"
$cv.Range("E5").Value = "Andisols"
$cv.Range("E6").Value = "Aridisols"
$cv.Range("E7").Value = "Entisols"
$cv.Range("E8").Value = "Gelisols"
$cv.Range("E9").Value = "Histosols"
$cv.Range("E10").Value = "Inceptisols"
$cv.Range("E11").Value = "Mollisols"
$cv.Range("E12").Value = "Oxisols"
$cv.Range("E13").Value = "Spodosols"
$cv.Range("E14").Value = "Ultisols"
$cv.Range("E15").Value = "Vertisols"

# --- 3. Fill in the existing profile rows with the previously-assigned
#        soil order (Vertisols) and attach the dropdown list validation
#        to the whole column for future rows ---------------------------
$profile.Range("N4:N7").Value = "Vertisols"

$validationRange = $profile.Range("N4:N1048576")
$validationRange.Validation.Add(3, 1, 1, "='controlled vocabulary'!`$E`$4:`$E`$15")
